# Update the "想去人数" (interested count) values in column F across the
# relevant sheets, matching the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 26
$ws1.Range("F3").Value = 8187
$ws1.Range("F4").Value = 1919
$ws1.Range("F7").Value = 2077
$ws1.Range("F8").Value = 572
$ws1.Range("F11").Value = 54
$ws1.Range("F15").Value = 5
$ws1.Range("F16").Value = 8569
$ws1.Range("F21").Value = 1813
$ws1.Range("F26").Value = 39
$ws1.Range("F30").Value = 6
$ws1.Range("F33").Value = 2105
$ws1.Range("F35").Value = 489
$ws1.Range("F37").Value = 1
$ws1.Range("F39").Value = 190
$ws1.Range("F42").Value = 40
$ws1.Range("F44").Value = 35

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2326

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2326
$ws4.Range("F5").Value = 8187
$ws4.Range("F8").Value = 1919
$ws4.Range("F11").Value = 2077
$ws4.Range("F14").Value = 572
$ws4.Range("F18").Value = 54
$ws4.Range("F23").Value = 8569
$ws4.Range("F27").Value = 1813
$ws4.Range("F32").Value = 2105
$ws4.Range("F35").Value = 489
